# Commit of the day
# - Fix "ChatGPT 3.5. Turbo" -> "ChatGPT 3.5 Turbo" typo
# - Remove the duplicate "Sheet" (total_time) and duplicate "Only US-Result" sheets
# - Rename the surviving "Only US-Result sheet " tab to "Only US-Result"
# - Rename "TimeConsupNoAnn" to "Time Consup."

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Fix the typo in the surviving "Only US-Result sheet " worksheet.
$wsResult = $wb.Worksheets.Item("Only US-Result sheet ")
$wsResult.Range("A2").Value = "ChatGPT 3.5 Turbo"
$wsResult.Range("A3").Value = "ChatGPT 3.5 Turbo"

# Remove the redundant "Sheet" (time totals) worksheet.
$wb.Worksheets.Item("Sheet").Delete() | Out-Null

# Remove the duplicate "Only US-Result" worksheet (identical data, already
# had the corrected text) now that the first sheet carries the fix.
$wb.Worksheets.Item("Only US-Result").Delete() | Out-Null

# Rename the remaining sheets to match the cleaned-up naming scheme.
$wsResult.Name = "Only US-Result"
$wb.Worksheets.Item("TimeConsupNoAnn").Name = "Time Consup."
